# "BASIC lighting functions no CAN"
# The BRAKE_EN, BACKUP_EN and EXT_SW_EN pins move from plain GPIO outputs to
# TIM3 PWM alternate-function outputs, so their Type/Description/Notes text
# is updated to reflect the new AF_PWM usage.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 - BRAKE_EN (PB0): GPIO OUT -> AF_PWM, new TIM3 Ch3 description, PWM note
$ws.Range("C4").Value = "AF_PWM"
# set the "PWM push pull" note cells first (E5 without trailing space, then E4
# with trailing space) so the shared-string table fills in the same order as
# the authored workbook
$ws.Range("E5").Value = "PWM push pull"
$ws.Range("E4").Value = "PWM push pull "
$ws.Range("D4").Value = "Brake light enable pin. TIM3 Ch3 "

# Row 5 - BACKUP_EN (PB1): GPIO OUT -> AF_PWM, new TIM3 Ch4 description
$ws.Range("D5").Value = "Backup light enable pin. TIM3 Ch4"
$ws.Range("C5").Value = "AF_PWM"

# Row 10 - EXT_SW_EN (PA6): GPIO OUT -> AF_PWM, new TIM3 Ch1 (TURN EN) description
$ws.Range("C10").Value = "AF_PWM"
$ws.Range("D10").Value = "EXT Switch TIM3 Ch1  (TURN EN)"

# Move the active selection like the author's last save
$ws.Range("A22").Select()
